$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.950.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.817.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4661"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3703"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07381"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8714"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.852.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.361"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07085"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.513"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008721"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.98%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.959.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.354"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.132.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.896"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.84"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.183"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.313"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08917"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7631"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.164"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.67%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.485"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.74%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.932"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.095"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01960"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05260"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.55%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5363"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.59%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.404"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.250"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.913"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1663"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.512"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4966"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.677"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06288"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.49%  "
